# Commit 0.0.12: Throw IllegalArgumentException instead of NullPointerException,
# on delegate issue errors. Adds a new error-message row (XML2SOURCE_FILE.ERR007)
# to the "ja" resource-bundle worksheet, just after the ERR006 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ja")

# Insert a brand new row at position 55; this pushes the former rows 55-87
# down to 56-88 and keeps all their formulas/values/styles intact.
$ws.Rows(55).Insert()

# The freshly inserted row 55 has no formatting of its own yet - copy the
# number formats/styles/borders from row 56 (the row that used to be the
# original row 55) so it matches its neighbours exactly.
$ws.Range("A56:G56").Copy()
$ws.Range("A55:G55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row counter formula, consistent with the "+1 from the row above" pattern
# used throughout the table.
$ws.Range("A55").Formula = "=A54+1"

# New error-message key/text pair for the new row.
$ws.Range("B55").Value = "XML2SOURCE_FILE.ERR007"
$ws.Range("C55").Value = "クラス名[{0}]の委譲フィールド[{1}]の型名が指定されていません。"

$excel.Calculate()

# Restore the view's selection to roughly where editing left off.
$ws.Range("C57").Select()
